# Apply the strategy-label and row-swap edits described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the strategy name text (shared strings) in column B, rows 2-8.
$ws.Range("B2").Value = "DB Search"
$ws.Range("B3").Value = "SB Search (BS*FS)"
$ws.Range("B4").Value = "DB Search + BS*FS"
$ws.Range("B5").Value = "Scopus + BS*FS"
$ws.Range("B6").Value = "Scopus + BS||FS"
$ws.Range("B7").Value = "Scopus + BS+FS"
$ws.Range("B8").Value = "Scopus + FS+BS"

# 2) Swap the numeric data (columns C:K) between row 3 and row 4, keeping
#    columns A and B (index + label) as they are.
$row3 = @()
$row4 = @()
for ($col = 3; $col -le 11; $col++) {
    $row3 += ,$ws.Cells.Item(3, $col).Value()
    $row4 += ,$ws.Cells.Item(4, $col).Value()
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(3, $col).Value = $row4[$i]
    $ws.Cells.Item(4, $col).Value = $row3[$i]
}
